$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,10
$data[0,0] = -19.51024841918096
$data[0,1] = 2.454197159004916
$data[0,2] = -19.51024841918096
$data[0,3] = -19.51024841918096
$data[0,4] = -19.51024841918096
$data[0,5] = -19.51024841918096
$data[0,6] = -19.51024841918096
$data[0,7] = -19.51024841918096
$data[0,8] = -19.51024841918096
$data[0,9] = -19.51024841918096
$data[1,0] = -19.51024841918096
$data[1,1] = -19.51024841918096
$data[1,2] = -19.51024841918096
$data[1,3] = -19.51024841918096
$data[1,4] = -19.51024841918096
$data[1,5] = -19.51024841918096
$data[1,6] = -19.51024841918096
$data[1,7] = 2.595279372573707
$data[1,8] = -19.51024841918096
$data[1,9] = -19.51024841918096
$data[2,0] = -19.51024841918096
$data[2,1] = 2.14794257535352
$data[2,2] = 2.823758834661553
$data[2,3] = -19.51024841918096
$data[2,4] = 2.600468410951405
$data[2,5] = -19.51024841918096
$data[2,6] = -19.51024841918096
$data[2,7] = -19.51024841918096
$data[2,8] = 2.067750579835688
$data[2,9] = -19.51024841918096
$data[3,0] = -19.51024841918096
$data[3,1] = 1.040574443817643
$data[3,2] = -19.51024841918096
$data[3,3] = -19.51024841918096
$data[3,4] = -19.51024841918096
$data[3,5] = 2.079847448550587
$data[3,6] = -19.51024841918096
$data[3,7] = -19.51024841918096
$data[3,8] = -19.51024841918096
$data[3,9] = -19.51024841918096
$data[4,0] = -19.51024841918096
$data[4,1] = -19.51024841918096
$data[4,2] = -19.51024841918096
$data[4,3] = -19.51024841918096
$data[4,4] = -19.51024841918096
$data[4,5] = -19.51024841918096
$data[4,6] = -19.51024841918096
$data[4,7] = -19.51024841918096
$data[4,8] = -19.51024841918096
$data[4,9] = -19.51024841918096
$data[5,0] = 2.969352516979972
$data[5,1] = -19.51024841918096
$data[5,2] = -19.51024841918096
$data[5,3] = -19.51024841918096
$data[5,4] = -19.51024841918096
$data[5,5] = -19.51024841918096
$data[5,6] = -19.51024841918096
$data[5,7] = -19.51024841918096
$data[5,8] = -19.51024841918096
$data[5,9] = -19.51024841918096
$data[6,0] = -19.51024841918096
$data[6,1] = -19.51024841918096
$data[6,2] = -19.51024841918096
$data[6,3] = 2.844683447092955
$data[6,4] = -19.51024841918096
$data[6,5] = -19.51024841918096
$data[6,6] = -19.51024841918096
$data[6,7] = -19.51024841918096
$data[6,8] = -19.51024841918096
$data[6,9] = -19.51024841918096
$data[7,0] = 3.605035476787425
$data[7,1] = -19.51024841918096
$data[7,2] = -19.51024841918096
$data[7,3] = -19.51024841918096
$data[7,4] = -19.51024841918096
$data[7,5] = -19.51024841918096
$data[7,6] = -19.51024841918096
$data[7,7] = -19.51024841918096
$data[7,8] = -19.51024841918096
$data[7,9] = -19.51024841918096
$data[8,0] = -19.51024841918096
$data[8,1] = -19.51024841918096
$data[8,2] = -19.51024841918096
$data[8,3] = -19.51024841918096
$data[8,4] = -19.51024841918096
$data[8,5] = -19.51024841918096
$data[8,6] = -19.51024841918096
$data[8,7] = 1.251430223797702
$data[8,8] = -19.51024841918096
$data[8,9] = 2.336091687435886
$data[9,0] = -19.51024841918096
$data[9,1] = -19.51024841918096
$data[9,2] = -19.51024841918096
$data[9,3] = 2.08086161302186
$data[9,4] = -19.51024841918096
$data[9,5] = 2.639134341772543
$data[9,6] = -19.51024841918096
$data[9,7] = -19.51024841918096
$data[9,8] = -19.51024841918096
$data[9,9] = 1.326040069668829
$data[10,0] = -19.51024841918096
$data[10,1] = -19.51024841918096
$data[10,2] = -19.51024841918096
$data[10,3] = -19.51024841918096
$data[10,4] = -19.51024841918096
$data[10,5] = -19.51024841918096
$data[10,6] = -19.51024841918096
$data[10,7] = -19.51024841918096
$data[10,8] = -19.51024841918096
$data[10,9] = -19.51024841918096
$data[11,0] = -19.51024841918096
$data[11,1] = -19.51024841918096
$data[11,2] = -19.51024841918096
$data[11,3] = 1.721948217178591
$data[11,4] = -19.51024841918096
$data[11,5] = -19.51024841918096
$data[11,6] = -19.51024841918096
$data[11,7] = -19.51024841918096
$data[11,8] = 2.341049602809201
$data[11,9] = 1.632557299268604
$data[12,0] = -19.51024841918096
$data[12,1] = -19.51024841918096
$data[12,2] = 1.70493439122561
$data[12,3] = -19.51024841918096
$data[12,4] = -19.51024841918096
$data[12,5] = -19.51024841918096
$data[12,6] = -19.51024841918096
$data[12,7] = -19.51024841918096
$data[12,8] = -19.51024841918096
$data[12,9] = 2.022113689573496
$data[13,0] = -19.51024841918096
$data[13,1] = -19.51024841918096
$data[13,2] = -0.191756361010128
$data[13,3] = -19.51024841918096
$data[13,4] = -19.51024841918096
$data[13,5] = -19.51024841918096
$data[13,6] = -19.51024841918096
$data[13,7] = -19.51024841918096
$data[13,8] = -19.51024841918096
$data[13,9] = -19.51024841918096
$data[14,0] = -19.51024841918096
$data[14,1] = -19.51024841918096
$data[14,2] = -19.51024841918096
$data[14,3] = -19.51024841918096
$data[14,4] = -19.51024841918096
$data[14,5] = -19.51024841918096
$data[14,6] = -19.51024841918096
$data[14,7] = -19.51024841918096
$data[14,8] = 2.23352302916645
$data[14,9] = -19.51024841918096
$data[15,0] = -19.51024841918096
$data[15,1] = 0.8769285493812802
$data[15,2] = 0.141442734592978
$data[15,3] = -19.51024841918096
$data[15,4] = -19.51024841918096
$data[15,5] = -19.51024841918096
$data[15,6] = -19.51024841918096
$data[15,7] = 0.9550638862524646
$data[15,8] = 1.314481185096337
$data[15,9] = -19.51024841918096
$data[16,0] = -19.51024841918096
$data[16,1] = -19.51024841918096
$data[16,2] = -19.51024841918096
$data[16,3] = -19.51024841918096
$data[16,4] = -19.51024841918096
$data[16,5] = -19.51024841918096
$data[16,6] = -19.51024841918096
$data[16,7] = 1.046928242094808
$data[16,8] = 1.828227378308875
$data[16,9] = -19.51024841918096
$data[17,0] = -19.51024841918096
$data[17,1] = -19.51024841918096
$data[17,2] = 1.756174991405372
$data[17,3] = -19.51024841918096
$data[17,4] = -19.51024841918096
$data[17,5] = -19.51024841918096
$data[17,6] = -19.51024841918096
$data[17,7] = 1.988963429052576
$data[17,8] = -19.51024841918096
$data[17,9] = -19.51024841918096
$data[18,0] = -19.51024841918096
$data[18,1] = 1.58966758399513
$data[18,2] = 2.10542117596129
$data[18,3] = -19.51024841918096
$data[18,4] = 3.800655677633801
$data[18,5] = -19.51024841918096
$data[18,6] = -19.51024841918096
$data[18,7] = 1.84873377196476
$data[18,8] = -19.51024841918096
$data[18,9] = 2.400841058934657
$data[19,0] = -19.51024841918096
$data[19,1] = 1.671130139030324
$data[19,2] = -19.51024841918096
$data[19,3] = 2.402471739000227
$data[19,4] = -19.51024841918096
$data[19,5] = 3.254408786962401
$data[19,6] = 4.321926259497099
$data[19,7] = -19.51024841918096
$data[19,8] = -19.51024841918096
$data[19,9] = -19.51024841918096

$ws.Range("B2:K21").Value = $data
